$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = 'Última actualización: 15:58:05'

$ws1.Range("A3").Value = 'Total filas: 239'

$ws1.Range("A33").Value = '06:17:28'
$ws1.Range("B33").Value = '08:00'
$ws1.Range("C33").Value = '17_ROMERO'
$ws1.Range("D33").Value = 103
$ws1.Range("E33").Value = 'LP1912'

$ws1.Range("A34").Value = '06:46:50'
$ws1.Range("B34").Value = '08:00'
$ws1.Range("C34").Value = '16_SANTA ANA'
$ws1.Range("D34").Value = 74
$ws1.Range("E34").Value = 'LP1912'

$ws1.Range("A64").Value = '08:27:16'
$ws1.Range("B64").Value = '09:39'
$ws1.Range("C64").Value = '23_HERNANDEZ'
$ws1.Range("D64").Value = 72
$ws1.Range("E64").Value = 'LP1912'

$ws1.Range("A65").Value = '07:50:33'
$ws1.Range("B65").Value = '09:39'
$ws1.Range("C65").Value = '15_ABASTO'
$ws1.Range("D65").Value = 109
$ws1.Range("E65").Value = 'LP1912'

$ws1.Range("A106").Value = '11:47:17'
$ws1.Range("B106").Value = '11:52'
$ws1.Range("C106").Value = '23_HERNANDEZ'
$ws1.Range("D106").Value = 5
$ws1.Range("E106").Value = 'LP1912'

$ws1.Range("A107").Value = '11:52:01'
$ws1.Range("B107").Value = '11:52'
$ws1.Range("C107").Value = '15X38_ABASTO'
$ws1.Range("D107").Value = 0
$ws1.Range("E107").Value = 'LP1912'

$ws1.Range("A108").Value = '10:05:51'
$ws1.Range("B108").Value = '11:52'
$ws1.Range("C108").Value = '225_GOMEZ'
$ws1.Range("D108").Value = 107
$ws1.Range("E108").Value = 'LP1912'

$ws1.Range("A120").Value = '10:37:52'
$ws1.Range("B120").Value = '12:10'
$ws1.Range("C120").Value = '16_P MOR-SANTA ANA'
$ws1.Range("D120").Value = 93
$ws1.Range("E120").Value = 'LP1912'

$ws1.Range("A121").Value = '10:37:52'
$ws1.Range("B121").Value = '12:10'
$ws1.Range("C121").Value = '15_ABASTO'
$ws1.Range("D121").Value = 93
$ws1.Range("E121").Value = 'LP1912'

$ws1.Range("A139").Value = '11:34:59'
$ws1.Range("B139").Value = '12:47'
$ws1.Range("C139").Value = '15X38_ABASTO'
$ws1.Range("D139").Value = 73
$ws1.Range("E139").Value = 'LP1912'

$ws1.Range("A141").Value = '11:34:59'
$ws1.Range("B141").Value = '12:47'
$ws1.Range("C141").Value = '14_ABASTO'
$ws1.Range("D141").Value = 73
$ws1.Range("E141").Value = 'LP1912'

$ws1.Range("A169").Value = '13:14:29'
$ws1.Range("B169").Value = '14:02'
$ws1.Range("C169").Value = '16_SANTA ANA'
$ws1.Range("D169").Value = 48
$ws1.Range("E169").Value = 'LP1912'

$ws1.Range("A170").Value = '12:33:21'
$ws1.Range("B170").Value = '14:02'
$ws1.Range("C170").Value = '10_OLMOS'
$ws1.Range("D170").Value = 89
$ws1.Range("E170").Value = 'LP1912'

$ws1.Range("A215").Value = '15:58:05'
$ws1.Range("B215").Value = '16:00'
$ws1.Range("C215").Value = '10_OLMOS'
$ws1.Range("D215").Value = 2
$ws1.Range("E215").Value = 'LP1912'

$ws1.Range("A216").Value = '14:12:26'
$ws1.Range("B216").Value = '16:05'
$ws1.Range("C216").Value = '14_ABASTO'
$ws1.Range("D216").Value = 113
$ws1.Range("E216").Value = 'LP1912'

$ws1.Range("A217").Value = '15:17:33'
$ws1.Range("B217").Value = '16:05'
$ws1.Range("C217").Value = '16_SANTA ANA'
$ws1.Range("D217").Value = 48
$ws1.Range("E217").Value = 'LP1912'

$ws1.Range("A218").Value = '14:53:07'
$ws1.Range("B218").Value = '16:06'
$ws1.Range("C218").Value = '14_ABASTO'
$ws1.Range("D218").Value = 73
$ws1.Range("E218").Value = 'LP1912'

$ws1.Range("A219").Value = '14:53:07'
$ws1.Range("B219").Value = '16:13'
$ws1.Range("C219").Value = '16_SANTA ANA'
$ws1.Range("D219").Value = 80
$ws1.Range("E219").Value = 'LP1912'

$ws1.Range("A220").Value = '14:32:38'
$ws1.Range("B220").Value = '16:14'
$ws1.Range("C220").Value = '17_ROMERO'
$ws1.Range("D220").Value = 102
$ws1.Range("E220").Value = 'LP1912'

$ws1.Range("A221").Value = '14:32:38'
$ws1.Range("B221").Value = '16:17'
$ws1.Range("C221").Value = '10_OLMOS'
$ws1.Range("D221").Value = 105
$ws1.Range("E221").Value = 'LP1912'

$ws1.Range("A222").Value = '14:53:07'
$ws1.Range("B222").Value = '16:20'
$ws1.Range("C222").Value = '23_HERNANDEZ'
$ws1.Range("D222").Value = 87
$ws1.Range("E222").Value = 'LP1912'

$ws1.Range("A223").Value = '14:32:38'
$ws1.Range("B223").Value = '16:21'
$ws1.Range("C223").Value = '23_HERNANDEZ'
$ws1.Range("D223").Value = 109
$ws1.Range("E223").Value = 'LP1912'

$ws1.Range("A224").Value = '15:46:07'
$ws1.Range("B224").Value = '16:22'
$ws1.Range("C224").Value = '23_HERNANDEZ'
$ws1.Range("D224").Value = 36
$ws1.Range("E224").Value = 'LP1912'

$ws1.Range("A225").Value = '15:17:33'
$ws1.Range("B225").Value = '16:30'
$ws1.Range("C225").Value = '16_SANTA ANA'
$ws1.Range("D225").Value = 73
$ws1.Range("E225").Value = 'LP1912'

$ws1.Range("A226").Value = '15:46:07'
$ws1.Range("B226").Value = '16:30'
$ws1.Range("C226").Value = '14_ABASTO'
$ws1.Range("D226").Value = 44
$ws1.Range("E226").Value = 'LP1912'

$ws1.Range("A227").Value = '14:45:56'
$ws1.Range("B227").Value = '16:33'
$ws1.Range("C227").Value = '83_ALUAR'
$ws1.Range("D227").Value = 108
$ws1.Range("E227").Value = 'LP1912'

$ws1.Range("A228").Value = '14:53:07'
$ws1.Range("B228").Value = '16:34'
$ws1.Range("C228").Value = '83_ALUAR'
$ws1.Range("D228").Value = 101
$ws1.Range("E228").Value = 'LP1912'

$ws1.Range("A229").Value = '14:45:56'
$ws1.Range("B229").Value = '16:40'
$ws1.Range("C229").Value = '225_GOMEZ'
$ws1.Range("D229").Value = 115
$ws1.Range("E229").Value = 'LP1912'

$ws1.Range("A230").Value = '14:53:07'
$ws1.Range("B230").Value = '16:41'
$ws1.Range("C230").Value = '225_GOMEZ'
$ws1.Range("D230").Value = 108
$ws1.Range("E230").Value = 'LP1912'

$ws1.Range("A231").Value = '14:53:07'
$ws1.Range("B231").Value = '16:46'
$ws1.Range("C231").Value = '17_ROMERO'
$ws1.Range("D231").Value = 113
$ws1.Range("E231").Value = 'LP1912'

$ws1.Range("A232").Value = '15:17:33'
$ws1.Range("B232").Value = '16:53'
$ws1.Range("C232").Value = '11_ETCHEVERRY'
$ws1.Range("D232").Value = 96
$ws1.Range("E232").Value = 'LP1912'

$ws1.Range("A233").Value = '15:46:07'
$ws1.Range("B233").Value = '16:54'
$ws1.Range("C233").Value = '11_ETCHEVERRY'
$ws1.Range("D233").Value = 68
$ws1.Range("E233").Value = 'LP1912'

$ws1.Range("A234").Value = '15:17:33'
$ws1.Range("B234").Value = '16:58'
$ws1.Range("C234").Value = '15_ABASTO'
$ws1.Range("D234").Value = 101
$ws1.Range("E234").Value = 'LP1912'

$ws1.Range("A235").Value = '15:17:33'
$ws1.Range("B235").Value = '17:07'
$ws1.Range("C235").Value = '16_P MOR-SANTA ANA'
$ws1.Range("D235").Value = 110
$ws1.Range("E235").Value = 'LP1912'

$ws1.Range("A236").Value = '15:46:07'
$ws1.Range("B236").Value = '17:09'
$ws1.Range("C236").Value = '23_HERNANDEZ'
$ws1.Range("D236").Value = 83
$ws1.Range("E236").Value = 'LP1912'

$ws1.Range("A237").Value = '15:17:33'
$ws1.Range("B237").Value = '17:10'
$ws1.Range("C237").Value = '215C_EL PATO'
$ws1.Range("D237").Value = 113
$ws1.Range("E237").Value = 'LP1912'

$ws1.Range("A238").Value = '15:46:07'
$ws1.Range("B238").Value = '17:21'
$ws1.Range("C238").Value = '15X38_ABASTO'
$ws1.Range("D238").Value = 95
$ws1.Range("E238").Value = 'LP1912'

$ws1.Range("A239").Value = '15:46:07'
$ws1.Range("B239").Value = '17:34'
$ws1.Range("C239").Value = '17_ROMERO'
$ws1.Range("D239").Value = 108
$ws1.Range("E239").Value = 'LP1912'

$ws1.Range("A240").Value = '15:58:05'
$ws1.Range("B240").Value = '17:36'
$ws1.Range("C240").Value = '27_EL RETIRO'
$ws1.Range("D240").Value = 98
$ws1.Range("E240").Value = 'LP1912'

$ws1.Range("A241").Value = '15:46:07'
$ws1.Range("B241").Value = '17:37'
$ws1.Range("C241").Value = '27_EL RETIRO'
$ws1.Range("D241").Value = 111
$ws1.Range("E241").Value = 'LP1912'

$ws1.Range("A242").Value = '15:46:07'
$ws1.Range("B242").Value = '17:39'
$ws1.Range("C242").Value = '215B_EL PATO'
$ws1.Range("D242").Value = 113
$ws1.Range("E242").Value = 'LP1912'

$ws1.Range("A243").Value = '15:46:07'
$ws1.Range("B243").Value = '17:42'
$ws1.Range("C243").Value = '215_EL PELIGRO'
$ws1.Range("D243").Value = 116
$ws1.Range("E243").Value = 'LP1912'

$ws1.Range("A244").Value = '15:58:05'
$ws1.Range("B244").Value = '17:46'
$ws1.Range("C244").Value = '215_EL PELIGRO'
$ws1.Range("D244").Value = 108
$ws1.Range("E244").Value = 'LP1912'

# --- Sheet: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = 'Última actualización: 15:58:05'

$ws2.Range("A3").Value = 'Total filas: 38'

$ws2.Range("A43").Value = '15:58:05'
$ws2.Range("B43").Value = '17:46'
$ws2.Range("C43").Value = '215_EL PELIGRO'
$ws2.Range("D43").Value = 108
$ws2.Range("E43").Value = 'LP1912'
